$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Gamma2F"

# Tiny value corrections on row 13 (precision tweaks)
$ws.Range("D13").Value2 = 0.9949169638402917
$ws.Range("J13").Value2 = 0.9949169638402917
$ws.Range("K13").Value2 = 0.9946409653551035

# New row 16 of data (reuse the "HexGrid-60degTilt5degRes" label from row 15)
$ws.Range("A16").Value2 = 14
$ws.Range("B16").Value2 = $ws.Range("B15").Value2

$ws.Range("C16").Value2 = 1.18355117195812
$ws.Range("D16").Value2 = 0.5960752793557353
$ws.Range("E16").Value2 = 1.047915988402299
$ws.Range("F16").Value2 = 1.18355117195812
$ws.Range("G16").Value2 = 0.7931908076294506
$ws.Range("H16").Value2 = 1.124909215325272
$ws.Range("I16").Value2 = 1.092202312283082
$ws.Range("J16").Value2 = 0.5960752793557353
$ws.Range("K16").Value2 = 0.8219956338790171
$ws.Range("L16").Value2 = 1.002773402918569
$ws.Range("M16").Value2 = 0.9729741291589932

# Copy style from A15 (bold/bordered/centered) to A16
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
